# Add a "Total" column header (F1) to each worksheet, matching the
# commit "updated data to include total plan".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PBO")
$ws1.Range("F1").Value = "Total"
$ws1.Range("F1").Select() | Out-Null

$ws2 = $wb.Worksheets.Item("Service Cost")
$ws2.Range("F1").Value = "Total"
$ws2.Range("F1").Select() | Out-Null

$ws3 = $wb.Worksheets.Item("PVFB")
$ws3.Range("F1").Value = "Total"
$ws3.Select() | Out-Null
$ws3.Range("E3").Select() | Out-Null
